# Generate Report for Handback
#
# 1. Update the "Status" column text everywhere it says "Ready for handoff"
#    to "Handed back: in sync with en-US" (Overview!B2:B3/C2:C3,
#    zh-cn!C2:C3, de-de!C2:C3).
# 2. Fill in "Latest Target File" (F) / "Latest Handback File" (G) columns
#    on the zh-cn and de-de sheets for rows 2-3 with hyperlinked file
#    names (mirrors the existing Source File / Latest Handoff File links).
# 3. Stamp the "Latest Handback DateTime" (H) column with real timestamps
#    (was the zero-date placeholder "0001-01-01 00:00:00").

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# --- 1. Status text swap on every sheet that shows it -----------------
$ovw = $wb.Worksheets.Item("Overview")
foreach ($addr in @("B2", "C2", "B3", "C3")) {
    $cell = $ovw.Range($addr)
    if ($cell.Value() -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3")) {
        $cell = $ws.Range($addr)
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- 2 & 3. Per-locale target/handback columns + handback timestamp ---
# locale -> @{ MdHost = base md hyperlink target prefix
#               XlfHost = locale-specific xlf hyperlink target prefix
#               HandbackTime2 = H2 new value
#               HandbackTime3 = H3 new value (same timestamp as row2) }
$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/c769b5c4d753f3d4eabb407af6c66439fc25749c/e2e/"

$locales = @{
    "zh-cn" = @{
        XlfHost = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b400d5ae2bd3bca6d3aacb54c20174ecd2cfdaf3/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/"
        Ext = "zh-cn.xlf"
        HandbackTime = "2016-03-23 06:53:51"
    }
    "de-de" = @{
        XlfHost = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/408a65326922755a8a594906d7c12f96ebc2cfb5/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/"
        Ext = "de-de.xlf"
        HandbackTime = "2016-03-23 06:54:04"
    }
}

$rows = @(
    @{ Row = 2; Uuid = "8f63d5ce-5ec1-4340-9423-1e8ba24d4048"; Hash = "919bc775bc0255afb34d99f0f3b7eac8ccd7776a" },
    @{ Row = 3; Uuid = "aeb33cbb-bc4d-40cb-ba3e-28416ed52bfb"; Hash = "86e58d6e1b66009a1da487a5ebb7ef60f3114c98" }
)

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $locales[$sheetName]

    foreach ($r in $rows) {
        $row = $r.Row
        $uuid = $r.Uuid
        $hash = $r.Hash

        $mdName = "$uuid.md"
        $xlfName = "$uuid.$hash.$($info.Ext)"

        $mdUrl = "$mdBase$mdName"
        $xlfUrl = "$($info.XlfHost)$xlfName"

        # Latest Target File (F) - same file the source points to
        $fCell = $ws.Range("F$row")
        $fCell.Value = $mdName
        $ws.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdName) | Out-Null

        # Latest Handback File (G) - same xlf the handoff points to
        $gCell = $ws.Range("G$row")
        $gCell.Value = $xlfName
        $ws.Hyperlinks.Add($gCell, $xlfUrl, "", "", $xlfName) | Out-Null

        # Latest Handback DateTime (H)
        $ws.Range("H$row").Value = $info.HandbackTime
    }
}
